# Journal entry method refactored
# - Set B8 on "Regression Suite" sheet to "Journal page"
# - Adjust row 3 height to 102.75 (custom height)
# - Select B8 as the active cell

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Regression Suite")

# Set the value of B8 (row 8, col 2) to "Journal page"
$ws.Range("B8").Value = "Journal page"

# Adjust row 3's height to a custom height of 102.75
$ws.Rows.Item(3).RowHeight = 102.75

# Activate the sheet and select B8 so it becomes the active cell/selection
$ws.Activate()
$ws.Range("B8").Select()
